# STATS.xlsx update: fill in second-pass (station-by-station then station-by-line)
# data for several lines/trams, fix "RHOME" -> "RHONE" typo, and tidy up a couple
# of stray leftover values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 (LIGNE A): add the missing "time as day fraction" helper formula ---
$ws.Range("I4").Formula = "=H4/86400"
$ws.Range("I4").NumberFormat = "hh:\ mm:\ ss"

# --- Row 7 (LIGNE D): add helper + success-rate formulas (H7/J7/K7 still blank) ---
$ws.Range("I7").Formula = "=H7/86400"
$ws.Range("I7").NumberFormat = "hh:\ mm:\ ss"
$ws.Range("M7").Formula = "=J7/L7"
$ws.Range("M7").NumberFormat = "0.00%"

# --- Row 8 (METRO COMMUNS): same treatment ---
$ws.Range("I8").Formula = "=H8/86400"
$ws.Range("I8").NumberFormat = "hh:\ mm:\ ss"
$ws.Range("M8").Formula = "=J8/L8"
$ws.Range("M8").NumberFormat = "0.00%"

# --- Row 9 (METRO ABC COMMUNS): fill in the actual second-pass figures ---
$ws.Range("F9").Formula = "=D9+E9"
$ws.Range("H9").Value = 168
$ws.Range("I9").Formula = "=H9/86400"
$ws.Range("I9").NumberFormat = "hh:\ mm:\ ss"
$ws.Range("J9").Value = 3040
$ws.Range("K9").Value = 338
$ws.Range("M9").Formula = "=J9/L9"
$ws.Range("M9").NumberFormat = "0.00%"

# --- Row 10 (METRO C COMMUNS) ---
$ws.Range("I10").Formula = "=H10/86400"
$ws.Range("I10").NumberFormat = "hh:\ mm:\ ss"
$ws.Range("M10").Formula = "=J10/L10"
$ws.Range("M10").NumberFormat = "0.00%"

# --- Row 11 (LIGNE C INTERSTATIONS) ---
$ws.Range("I11").Formula = "=H11/86400"
$ws.Range("I11").NumberFormat = "hh:\ mm:\ ss"
$ws.Range("M11").Formula = "=J11/L11"
$ws.Range("M11").NumberFormat = "0.00%"

# --- Row 12 (T1): corrected first-pass counts + new second-pass figures ---
$ws.Range("D12").Value = 4612
$ws.Range("E12").Value = 361
$ws.Range("H12").Value = 260
$ws.Range("I12").Formula = "=H12/86400"
$ws.Range("I12").NumberFormat = "hh:\ mm:\ ss"
$ws.Range("J12").Value = 4675
$ws.Range("K12").Value = 301
$ws.Range("M12").Formula = "=J12/L12"
$ws.Range("M12").NumberFormat = "0.00%"

# --- Row 13 (T2): new second-pass figures ---
$ws.Range("H13").Value = 914
$ws.Range("I13").Formula = "=H13/86400"
$ws.Range("I13").NumberFormat = "hh:\ mm:\ ss"
$ws.Range("J13").Value = 7223
$ws.Range("K13").Value = 140
$ws.Range("M13").Formula = "=J13/L13"
$ws.Range("M13").NumberFormat = "0.00%"

# --- Row 14 (T3): corrected first-pass counts + new second-pass figures ---
$ws.Range("D14").Value = 11000
$ws.Range("E14").Value = 802
$ws.Range("H14").Value = 483
$ws.Range("I14").Formula = "=H14/86400"
$ws.Range("I14").NumberFormat = "hh:\ mm:\ ss"
$ws.Range("J14").Value = 10134
$ws.Range("K14").Value = 674
$ws.Range("M14").Formula = "=J14/L14"
$ws.Range("M14").NumberFormat = "0.00%"

# --- Row 19: fix the "RHOME" typo to "RHONE" ---
$ws.Range("A19").Value = "RHONE EXPRESS COMMUNS"

# --- Row 28 (SERBER 500000): corrected second-pass counts ---
$ws.Range("J28").Value = 453
$ws.Range("K28").Value = 5

# --- Row 29 (SERBER 100000): corrected second-pass counts ---
$ws.Range("J29").Value = 4329
$ws.Range("K29").Value = 38

# --- Row 31: remove the stray leftover value ---
$ws.Range("J31").ClearContents()

# --- Update the active selection to reflect where the author ended up ---
$ws.Range("A19").Select()
